$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (C1:F1) and add new headers (G1:J1) ---
$ws.Cells.Item(1, 3).Value = "R2 (val)"
$ws.Cells.Item(1, 4).Value = "MAPE (val)"
$ws.Cells.Item(1, 5).Value = "RMSE (val)"
$ws.Cells.Item(1, 6).Value = "%RMSE (val)"
$ws.Cells.Item(1, 7).Value = "R2 (train)"
$ws.Cells.Item(1, 8).Value = "MAPE (train)"
$ws.Cells.Item(1, 9).Value = "RMSE (train)"
$ws.Cells.Item(1, 10).Value = "%RMSE (train)"

# Copy header style (bold, border, center) from F1 to the new header cells G1:J1
$ws.Range("F1").Copy()
$ws.Range("G1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update data rows 2-57 (columns C-J) ---
$data = @(
  @("0.32 ± 0.09", "2.50 ± 1.52", "49.82 ± 11.06", "65.08 ± 15.61", "0.91 ± 0.02", "1.46 ± 0.26", "18.52 ± 1.32", "24.10 ± 1.76"),
  @("0.31 ± 0.08", "2.51 ± 1.72", "50.11 ± 9.59", "65.44 ± 13.82", "0.90 ± 0.01", "1.03 ± 0.18", "19.73 ± 0.76", "25.69 ± 1.29"),
  @("0.11 ± 0.10", "3.95 ± 3.58", "57.06 ± 12.68", "74.58 ± 18.20", "0.19 ± 0.04", "3.84 ± 0.85", "55.36 ± 3.64", "72.01 ± 3.75"),
  @("0.32 ± 0.07", "2.60 ± 1.81", "49.64 ± 9.95", "64.86 ± 14.44", "0.82 ± 0.01", "1.43 ± 0.28", "26.26 ± 1.27", "34.18 ± 1.66"),
  @("0.29 ± 0.14", "2.14 ± 0.94", "50.49 ± 10.45", "66.00 ± 15.46", "0.93 ± 0.00", "0.75 ± 0.18", "16.07 ± 0.50", "20.92 ± 0.64"),
  @("0.16 ± 0.21", "3.88 ± 3.82", "53.86 ± 6.93", "70.33 ± 10.98", "0.41 ± 0.05", "2.92 ± 0.72", "47.01 ± 1.27", "61.17 ± 0.89"),
  @("0.31 ± 0.11", "1.74 ± 0.61", "49.70 ± 9.38", "64.87 ± 13.34", "0.83 ± 0.01", "0.93 ± 0.16", "25.18 ± 1.33", "32.78 ± 1.56"),
  @("0.07 ± 0.21", "3.82 ± 4.07", "58.36 ± 14.99", "76.36 ± 21.34", "0.18 ± 0.03", "3.57 ± 0.77", "55.66 ± 3.31", "72.42 ± 3.71"),
  @("0.34 ± 0.07", "2.65 ± 1.86", "49.24 ± 10.61", "64.34 ± 15.10", "0.69 ± 0.01", "1.74 ± 0.23", "33.95 ± 1.29", "44.19 ± 1.59"),
  @("0.08 ± 0.25", "3.46 ± 3.68", "56.50 ± 6.91", "73.76 ± 11.13", "0.48 ± 0.04", "2.43 ± 0.76", "44.19 ± 1.39", "57.51 ± 1.30"),
  @("0.33 ± 0.08", "2.38 ± 1.59", "49.42 ± 9.83", "64.53 ± 14.03", "0.76 ± 0.03", "1.57 ± 0.30", "30.29 ± 2.63", "39.38 ± 2.56"),
  @("0.30 ± 0.10", "2.40 ± 1.47", "50.03 ± 9.24", "65.33 ± 13.41", "1.00 ± 0.00", "0.16 ± 0.05", "3.15 ± 0.56", "4.11 ± 0.86"),
  @("0.24 ± 0.13", "2.84 ± 2.22", "52.02 ± 10.14", "67.96 ± 14.74", "0.67 ± 0.02", "1.82 ± 0.39", "31.28 ± 1.61", "40.71 ± 1.80"),
  @("0.35 ± 0.08", "2.36 ± 1.39", "48.63 ± 10.21", "63.51 ± 14.48", "0.84 ± 0.01", "1.38 ± 0.21", "24.37 ± 1.04", "31.72 ± 1.31"),
  @("0.30 ± 0.06", "2.39 ± 1.43", "50.36 ± 10.53", "65.80 ± 15.14", "0.93 ± 0.01", "1.21 ± 0.27", "16.24 ± 0.82", "21.15 ± 1.13"),
  @("0.31 ± 0.09", "2.64 ± 1.86", "49.84 ± 8.97", "65.07 ± 13.01", "0.90 ± 0.01", "1.02 ± 0.13", "19.55 ± 0.62", "25.45 ± 1.01"),
  @("0.12 ± 0.10", "3.72 ± 3.20", "56.91 ± 13.32", "74.41 ± 19.04", "0.20 ± 0.05", "3.74 ± 0.80", "54.99 ± 3.85", "71.52 ± 3.89"),
  @("0.28 ± 0.14", "2.48 ± 1.58", "50.87 ± 9.38", "66.44 ± 13.90", "0.90 ± 0.00", "1.12 ± 0.21", "19.03 ± 0.57", "24.77 ± 0.85"),
  @("0.28 ± 0.12", "2.27 ± 1.18", "50.91 ± 9.97", "66.51 ± 14.56", "0.93 ± 0.00", "1.01 ± 0.20", "15.76 ± 0.83", "20.50 ± 0.90"),
  @("0.17 ± 0.23", "3.41 ± 3.15", "53.46 ± 5.99", "69.79 ± 9.93", "0.44 ± 0.04", "2.55 ± 0.64", "45.99 ± 1.17", "59.85 ± 0.56"),
  @("0.33 ± 0.06", "1.78 ± 0.62", "49.37 ± 9.44", "64.47 ± 13.61", "0.85 ± 0.00", "0.96 ± 0.11", "24.13 ± 1.05", "31.40 ± 0.90"),
  @("0.06 ± 0.30", "4.26 ± 3.99", "58.38 ± 15.62", "76.39 ± 21.88", "0.19 ± 0.08", "4.00 ± 0.94", "55.38 ± 5.39", "71.99 ± 5.78"),
  @("0.32 ± 0.05", "2.81 ± 1.99", "49.86 ± 10.48", "65.16 ± 15.10", "0.68 ± 0.01", "1.86 ± 0.25", "34.65 ± 1.31", "45.10 ± 1.42"),
  @("0.13 ± 0.19", "2.90 ± 2.84", "55.34 ± 7.81", "72.31 ± 12.45", "0.63 ± 0.02", "1.54 ± 0.69", "37.55 ± 1.44", "48.86 ± 1.00"),
  @("0.30 ± 0.06", "2.62 ± 1.73", "50.29 ± 9.73", "65.68 ± 13.99", "0.76 ± 0.03", "1.61 ± 0.38", "29.81 ± 3.29", "38.73 ± 3.33"),
  @("0.28 ± 0.11", "2.11 ± 0.99", "50.88 ± 8.53", "66.46 ± 12.78", "0.99 ± 0.00", "0.24 ± 0.05", "4.42 ± 0.70", "5.77 ± 1.06"),
  @("0.24 ± 0.13", "2.78 ± 2.05", "52.21 ± 9.98", "68.21 ± 14.62", "0.70 ± 0.02", "1.74 ± 0.39", "29.79 ± 1.75", "38.76 ± 1.82"),
  @("0.34 ± 0.05", "2.34 ± 1.36", "48.82 ± 9.65", "63.78 ± 13.93", "0.87 ± 0.01", "1.27 ± 0.20", "22.01 ± 0.78", "28.65 ± 0.84"),
  @("0.71 ± 0.09", "0.81 ± 0.40", "32.72 ± 10.67", "42.78 ± 14.66", "0.98 ± 0.00", "0.54 ± 0.15", "8.95 ± 0.81", "11.65 ± 1.14"),
  @("0.71 ± 0.08", "0.67 ± 0.35", "32.49 ± 9.22", "42.40 ± 12.44", "0.96 ± 0.01", "0.24 ± 0.05", "12.33 ± 0.98", "16.06 ± 1.51"),
  @("0.68 ± 0.08", "0.68 ± 0.22", "34.34 ± 10.00", "44.86 ± 13.71", "0.74 ± 0.02", "0.74 ± 0.09", "31.57 ± 2.68", "41.10 ± 3.49"),
  @("0.70 ± 0.07", "0.92 ± 0.73", "33.48 ± 9.72", "43.79 ± 13.72", "0.97 ± 0.00", "0.40 ± 0.08", "11.23 ± 1.02", "14.61 ± 1.18"),
  @("0.68 ± 0.07", "0.95 ± 0.49", "34.01 ± 8.77", "44.42 ± 12.21", "1.00 ± 0.00", "0.00 ± 0.00", "0.01 ± 0.00", "0.02 ± 0.00"),
  @("0.18 ± 0.18", "4.02 ± 4.02", "53.50 ± 7.21", "69.84 ± 11.22", "0.39 ± 0.05", "3.13 ± 0.68", "47.94 ± 1.44", "62.39 ± 1.49"),
  @("0.69 ± 0.05", "0.80 ± 0.43", "33.77 ± 8.28", "44.12 ± 11.63", "0.91 ± 0.01", "0.49 ± 0.09", "18.22 ± 1.96", "23.72 ± 2.50"),
  @("0.63 ± 0.08", "1.36 ± 0.97", "36.81 ± 9.95", "48.16 ± 14.06", "0.69 ± 0.06", "1.06 ± 0.40", "33.93 ± 4.14", "44.17 ± 5.43"),
  @("0.68 ± 0.07", "1.25 ± 0.86", "32.44 ± 8.94", "42.38 ± 12.29", "0.86 ± 0.01", "0.84 ± 0.12", "22.72 ± 1.67", "29.59 ± 2.29"),
  @("0.34 ± 0.10", "1.60 ± 1.52", "49.21 ± 10.67", "64.29 ± 15.16", "0.80 ± 0.02", "0.51 ± 0.31", "27.31 ± 2.54", "35.49 ± 2.58"),
  @("0.72 ± 0.09", "0.65 ± 0.29", "31.91 ± 9.21", "41.61 ± 12.24", "0.92 ± 0.02", "0.44 ± 0.15", "17.22 ± 1.87", "22.45 ± 2.89"),
  @("0.68 ± 0.11", "0.68 ± 0.40", "33.98 ± 11.14", "44.41 ± 15.16", "0.99 ± 0.00", "0.23 ± 0.06", "5.85 ± 0.50", "7.63 ± 0.83"),
  @("0.62 ± 0.09", "1.20 ± 0.89", "36.55 ± 9.48", "47.76 ± 13.21", "0.85 ± 0.02", "0.72 ± 0.18", "19.77 ± 1.63", "25.74 ± 2.11"),
  @("0.72 ± 0.06", "0.90 ± 0.61", "31.91 ± 9.29", "41.71 ± 12.88", "0.95 ± 0.01", "0.55 ± 0.09", "13.72 ± 1.22", "17.87 ± 1.58"),
  @("0.70 ± 0.10", "0.82 ± 0.46", "33.42 ± 11.35", "43.69 ± 15.41", "0.98 ± 0.00", "0.57 ± 0.15", "9.09 ± 0.80", "11.84 ± 1.09"),
  @("0.69 ± 0.08", "0.67 ± 0.34", "33.64 ± 9.71", "43.93 ± 13.27", "0.96 ± 0.01", "0.23 ± 0.03", "12.83 ± 1.16", "16.71 ± 1.63"),
  @("0.67 ± 0.06", "0.72 ± 0.16", "34.98 ± 9.85", "45.74 ± 13.79", "0.74 ± 0.02", "0.77 ± 0.13", "31.10 ± 2.63", "40.49 ± 3.43"),
  @("0.71 ± 0.07", "0.78 ± 0.47", "32.58 ± 9.62", "42.58 ± 13.30", "1.00 ± 0.00", "0.16 ± 0.02", "3.60 ± 0.23", "4.69 ± 0.35"),
  @("0.67 ± 0.06", "0.94 ± 0.49", "35.06 ± 9.20", "45.83 ± 12.99", "1.00 ± 0.00", "0.00 ± 0.00", "0.01 ± 0.00", "0.01 ± 0.00"),
  @("0.14 ± 0.26", "3.04 ± 3.06", "54.17 ± 5.09", "70.68 ± 8.61", "0.55 ± 0.01", "1.90 ± 0.59", "41.05 ± 2.04", "53.41 ± 1.93"),
  @("0.69 ± 0.05", "0.75 ± 0.36", "33.88 ± 8.69", "44.27 ± 12.15", "0.92 ± 0.01", "0.44 ± 0.08", "17.51 ± 1.93", "22.79 ± 2.43"),
  @("0.53 ± 0.07", "0.91 ± 0.32", "41.53 ± 10.01", "54.21 ± 13.72", "0.62 ± 0.09", "0.96 ± 0.40", "37.61 ± 4.24", "48.89 ± 4.77"),
  @("0.69 ± 0.05", "1.44 ± 1.01", "33.70 ± 8.98", "44.03 ± 12.43", "0.87 ± 0.01", "0.92 ± 0.13", "22.40 ± 1.68", "29.16 ± 2.17"),
  @("0.17 ± 0.17", "2.70 ± 2.53", "53.85 ± 7.88", "70.37 ± 12.47", "0.66 ± 0.01", "1.34 ± 0.63", "36.03 ± 1.51", "46.88 ± 1.16"),
  @("0.70 ± 0.08", "0.65 ± 0.24", "32.83 ± 9.25", "42.84 ± 12.42", "0.92 ± 0.02", "0.47 ± 0.19", "17.43 ± 2.12", "22.68 ± 2.63"),
  @("0.66 ± 0.11", "0.66 ± 0.34", "35.21 ± 11.77", "46.08 ± 16.35", "1.00 ± 0.00", "0.08 ± 0.02", "2.21 ± 0.36", "2.88 ± 0.48"),
  @("0.59 ± 0.10", "1.17 ± 0.81", "37.90 ± 9.28", "49.52 ± 13.08", "0.85 ± 0.02", "0.65 ± 0.20", "19.24 ± 1.56", "25.04 ± 1.84"),
  @("0.72 ± 0.07", "0.93 ± 0.56", "32.09 ± 9.61", "41.94 ± 13.19", "0.96 ± 0.00", "0.51 ± 0.08", "11.68 ± 0.96", "15.20 ± 1.24")
)

$r = 2
foreach ($row in $data) {
    $c = 3
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
